$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 836.4666999999999
$ws.Range("I39").Value = 142.25
$ws.Range("J39").Value = 1088.909
$ws.Range("K39").Value = 426.75
$ws.Range("L39").Value = 3266.727
$ws.Range("M39").Value = -130.75
$ws.Range("N39").Value = -3858.727

$ws.Range("H113").Value = 2430.5
$ws.Range("I113").Value = 2144
$ws.Range("K113").Value = 2144
$ws.Range("M113").Value = 1110

$ws.Range("H137").Value = 1308.3235
$ws.Range("I137").Value = 940.1
$ws.Range("J137").Value = 1834.3572
$ws.Range("K137").Value = 2820.3
$ws.Range("L137").Value = 5503.071599999999
$ws.Range("M137").Value = -270.3000000000002
$ws.Range("N137").Value = -10603.0716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 792.3
$ws.Range("I2").Value = 726.375
$ws.Range("J2").Value = 1056
$ws.Range("K2").Value = 726.375
$ws.Range("L2").Value = 1056
$ws.Range("M2").Value = -613.375
$ws.Range("N2").Value = -1282

$ws.Range("H32").Value = 2159372.2
$ws.Range("I32").Value = 2538773
$ws.Range("J32").Value = 18469.072
$ws.Range("K32").Value = 2538773
$ws.Range("L32").Value = 18469.072
$ws.Range("M32").Value = -2538486
$ws.Range("N32").Value = -19043.072

$ws.Range("H45").Value = 1614.0869
$ws.Range("I45").Value = 1421.6154
$ws.Range("J45").Value = 1864.3
$ws.Range("K45").Value = 1421.6154
$ws.Range("L45").Value = 1864.3
$ws.Range("M45").Value = -1044.6154
$ws.Range("N45").Value = -2618.3

$ws.Range("H61").Value = 8337500
$ws.Range("I61").Value = 17545852
$ws.Range("J61").Value = 6133.2856
$ws.Range("K61").Value = 17545852
$ws.Range("L61").Value = 6133.2856
$ws.Range("M61").Value = -17545640
$ws.Range("N61").Value = -6557.2856

$ws.Range("H86").Value = 41695396
$ws.Range("I86").Value = 29500
$ws.Range("J86").Value = 43506956
$ws.Range("K86").Value = 29500
$ws.Range("L86").Value = 43506956
$ws.Range("N86").Value = -43509328
$ws.Range("M86").Value = -28314

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H89").Value = 41695396
$ws.Range("I89").Value = 29500
$ws.Range("J89").Value = 43506956
$ws.Range("K89").Value = 88500
$ws.Range("L89").Value = 130520868
$ws.Range("N89").Value = -130532724
$ws.Range("M89").Value = -82572

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H116").Value = 792.3
$ws.Range("I116").Value = 726.375
$ws.Range("J116").Value = 1056
$ws.Range("K116").Value = 726.375
$ws.Range("L116").Value = 1056
$ws.Range("M116").Value = 1567.625
$ws.Range("N116").Value = -5644

$ws.Range("H122").Value = 43119.793
$ws.Range("I122").Value = 54060.26
$ws.Range("J122").Value = 1546
$ws.Range("K122").Value = 162180.78
$ws.Range("L122").Value = 4638
$ws.Range("M122").Value = -159730.78
$ws.Range("N122").Value = -9538

$ws.Range("H123").Value = 59762.332
$ws.Range("J123").Value = 59762.332
$ws.Range("L123").Value = 59762.332
$ws.Range("N123").Value = -69562.33199999999

$ws.Range("H136").Value = 8337500
$ws.Range("I136").Value = 17545852
$ws.Range("J136").Value = 6133.2856
$ws.Range("K136").Value = 52637556
$ws.Range("L136").Value = 18399.8568
$ws.Range("M136").Value = -52635006
$ws.Range("N136").Value = -23499.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 792.3
$ws.Range("I3").Value = 726.375
$ws.Range("J3").Value = 1056
$ws.Range("K3").Value = 726.375
$ws.Range("L3").Value = 1056
$ws.Range("M3").Value = -612.375
$ws.Range("N3").Value = -1284

$ws.Range("H99").Value = 1833.3334
$ws.Range("I99").Value = 1800
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1800
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -302
$ws.Range("N99").Value = -4996

$ws.Range("H134").Value = 3620.3333
$ws.Range("I134").Value = 3787.4348
$ws.Range("J134").Value = 3071.2856
$ws.Range("K134").Value = 11362.3044
$ws.Range("L134").Value = 9213.856800000001
$ws.Range("M134").Value = -8827.304400000001
$ws.Range("N134").Value = -14283.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4811.725
$ws.Range("I31").Value = 1262.7297
$ws.Range("J31").Value = 7865.5117
$ws.Range("K31").Value = 1262.7297
$ws.Range("L31").Value = 7865.5117
$ws.Range("M31").Value = -967.7297000000001
$ws.Range("N31").Value = -8455.511699999999

$ws.Range("H34").Value = 4811.725
$ws.Range("I34").Value = 1262.7297
$ws.Range("J34").Value = 7865.5117
$ws.Range("K34").Value = 1262.7297
$ws.Range("L34").Value = 7865.5117
$ws.Range("M34").Value = -1060.7297
$ws.Range("N34").Value = -8269.511699999999

$ws.Range("H58").Value = 2037.5
$ws.Range("I58").Value = 1550
$ws.Range("K58").Value = 1550
$ws.Range("M58").Value = -1347

$ws.Range("H107").Value = 3125795.2
$ws.Range("I107").Value = 6944812.5
$ws.Range("J107").Value = 1144.8182
$ws.Range("K107").Value = 6944812.5
$ws.Range("L107").Value = 1144.8182
$ws.Range("M107").Value = -6942892.5
$ws.Range("N107").Value = -4984.8182

$ws.Range("H116").Value = 30000
$ws.Range("J116").Value = 30000
$ws.Range("L116").Value = 30000
$ws.Range("N116").Value = -39178

$ws.Range("H122").Value = 1771.4
$ws.Range("J122").Value = 2036.25
$ws.Range("L122").Value = 6108.75
$ws.Range("N122").Value = -11008.75

$ws.Range("H136").Value = 2037.5
$ws.Range("I136").Value = 1550
$ws.Range("K136").Value = 4650
$ws.Range("M136").Value = -2100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 631.93335
$ws.Range("J5").Value = 968.8182
$ws.Range("L5").Value = 2906.4546
$ws.Range("N5").Value = -3130.4546

$ws.Range("H135").Value = 631.93335
$ws.Range("J135").Value = 968.8182
$ws.Range("L135").Value = 8719.363800000001
$ws.Range("N135").Value = -13789.3638

$ws.Range("H137").Value = 28178.559
$ws.Range("I137").Value = 6394.4546
$ws.Range("J137").Value = 51000
$ws.Range("K137").Value = 19183.3638
$ws.Range("L137").Value = 153000
$ws.Range("M137").Value = -14083.3638
$ws.Range("N137").Value = -163200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 73333.336
$ws.Range("J68").Value = 97500
$ws.Range("L68").Value = 97500
$ws.Range("N68").Value = -99122

$ws.Range("H69").Value = 84000
$ws.Range("J69").Value = 84000
$ws.Range("L69").Value = 84000
$ws.Range("N69").Value = -85498

$ws.Range("H71").Value = 73333.336
$ws.Range("J71").Value = 97500
$ws.Range("L71").Value = 292500
$ws.Range("N71").Value = -300612

$ws.Range("H72").Value = 84000
$ws.Range("J72").Value = 84000
$ws.Range("L72").Value = 252000
$ws.Range("N72").Value = -259488

$ws.Range("H102").Value = 1494.3
$ws.Range("I102").Value = 1542.875
$ws.Range("K102").Value = 1542.875
$ws.Range("M102").Value = 79.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 98000
$ws.Range("J87").Value = 98000
$ws.Range("L87").Value = 98000
$ws.Range("N87").Value = -100246

$ws.Range("H90").Value = 98000
$ws.Range("J90").Value = 98000
$ws.Range("L90").Value = 294000
$ws.Range("N90").Value = -305232

$ws.Range("H100").Value = 71752.62
$ws.Range("I100").Value = 84160
$ws.Range("J100").Value = 3512
$ws.Range("K100").Value = 84160
$ws.Range("L100").Value = 3512
$ws.Range("M100").Value = -83619
$ws.Range("N100").Value = -4594

$ws.Range("H136").Value = 4275118
$ws.Range("I136").Value = 1204.4828
$ws.Range("K136").Value = 3613.4484
$ws.Range("M136").Value = -1063.4484

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4283.4287
$ws.Range("I81").Value = 5711.6665
$ws.Range("J81").Value = 3212.25
$ws.Range("K81").Value = 11423.333
$ws.Range("L81").Value = 6424.5
$ws.Range("M81").Value = -10362.333
$ws.Range("N81").Value = -8546.5

$ws.Range("H84").Value = 4283.4287
$ws.Range("I84").Value = 5711.6665
$ws.Range("J84").Value = 3212.25
$ws.Range("K84").Value = 57116.665
$ws.Range("L84").Value = 32122.5
$ws.Range("M84").Value = -51812.665
$ws.Range("N84").Value = -42730.5

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H132").Value = 4276547
$ws.Range("I132").Value = 3149.3157
$ws.Range("K132").Value = 9447.947100000001
$ws.Range("M132").Value = -6917.947100000001

$ws.Range("H136").Value = 2708.2322
$ws.Range("I136").Value = 2326.3684
$ws.Range("J136").Value = 3514.389
$ws.Range("K136").Value = 6979.1052
$ws.Range("L136").Value = 10543.167
$ws.Range("M136").Value = -4429.1052
$ws.Range("N136").Value = -15643.167
